$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every changed cell before writing, so that
# numeric-looking strings (prices, volume codes) stay stored as text
# -- matching the inlineStr/text representation used in the workbook.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "249.14"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.451"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05684"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.382"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8079"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "One"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.01167"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8OneONEBestin24h"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1471"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07710"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03164"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03027"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09266"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.552"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001656"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04721"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006351"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005034"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001042"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001501"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0003203"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.773"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.426"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1305"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04073"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006942"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003503"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1043"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007861"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005910"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005506"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6830"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.008934"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01011"
